# Add the "events" tab to the data workbook.
#
# Target end-state (per the commit diff):
#   - the old "none" summary sheet (wide layout: one column per event date)
#     is removed
#   - a new "events" sheet is added at the end, holding the same
#     information but reshaped into a long/tidy layout (one row per
#     ref/measure/date combination), plus a new "No more problems" (ref 101)
#     event series
#   - final sheet order: week, month, events

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- remove the old "none" sheet -----------------------------------------
$wb.Worksheets.Item("none").Delete() | Out-Null

# --- add the new "events" sheet at the end of the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "events"

$ws = $wb.Worksheets.Item("events")

# header row
$ws.Range("A1").Value = "ref"
$ws.Range("B1").Value = "measure_name"
$ws.Range("C1").Value = "comment"
$ws.Range("D1").Value = "event_date_or_datetime"

# long-format event data: one row per ref/measure/date
$comment = "you can put a comment here"
$dates = @(43890, 44511, 44662, 44834)

$row = 2
foreach ($item in @(
        @{ ref = 99; measure = "Problems" },
        @{ ref = 101; measure = "No more problems" }
    )) {
    $first = $true
    foreach ($d in $dates) {
        $ws.Cells.Item($row, 1).Value = $item.ref
        $ws.Cells.Item($row, 2).Value = $item.measure
        if ($first) {
            $ws.Cells.Item($row, 3).Value = $comment
            $first = $false
        }
        $ws.Cells.Item($row, 4).Value = $d
        $row++
    }
}

# apply the existing date number format (style) to the new date column,
# reusing the workbook's built-in date style rather than creating a new one
$dateSource = $wb.Worksheets.Item("week").Range("D1")
$dateSource.Copy() | Out-Null
$ws.Range("D2:D9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# column widths, matching the authored sheet as closely as this engine's
# character-width rounding allows
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 21.276041666666668
$ws.Columns.Item(3).ColumnWidth = 24.944010416666668
$ws.Columns.Item(4).ColumnWidth = 22.053385416666668

# --- reorder sheets: week, month, events ----------------------------------
# (the delete + append above already leaves the tabs in this order:
#  week, month, events)

# the "month" sheet is the active/selected tab in the authored workbook
$wb.Worksheets.Item("month").Activate()
